# Apply updates described by the commit:
# "update to published CDA FHIR logical model with patches #241"

$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.0.0-sd-202312-matchbox-patch -> 2.0.0-sd-202406-matchbox-patch
$meta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"

# Date: 2024-03-12T18:28:21+01:00 -> 2024-06-19T17:47:42+02:00
$meta.Range("B8").Value = "2024-06-19T17:47:42+02:00"

# Contact: "No display for ContactDetail" -> full HL7 contact string
$meta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Binding Value Set for row 5: v3-SetOperator -> CDASetOperator
$elements.Range("Z5").Value = "http://hl7.org/cda/stds/core/ValueSet/CDASetOperator"

# Column Z width change (49.5 -> 51.21484375) to fit the new (longer) text.
# Note: the ColumnWidth COM setter in this runtime snaps to a 1/6-character
# pixel grid (offset by 5/6), so the nearest representable value to the
# target 51.21484375 is produced by requesting 50.333333 (-> 51.16666667).
$elements.Columns.Item(26).ColumnWidth = 50.333333
